$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 533.44446
$ws.Range("I99").Value = 546.8333
$ws.Range("J99").Value = 506.66666
$ws.Range("K99").Value = 1640.4999
$ws.Range("L99").Value = 1519.99998
$ws.Range("M99").Value = -142.4999
$ws.Range("N99").Value = -4515.999980000001

$ws.Range("H103").Value = 2266.7856
$ws.Range("I103").Value = 1024
$ws.Range("J103").Value = 2763.9
$ws.Range("K103").Value = 3072
$ws.Range("L103").Value = 8291.700000000001
$ws.Range("M103").Value = -2486
$ws.Range("N103").Value = -9463.700000000001

$ws.Range("H116").Value = 7940.5
$ws.Range("I116").Value = 6464.6665
$ws.Range("J116").Value = 8826
$ws.Range("K116").Value = 6464.6665
$ws.Range("L116").Value = 8826
$ws.Range("M116").Value = -3022.6665
$ws.Range("N116").Value = -15710

$ws.Range("H131").Value = 4374.8096
$ws.Range("I131").Value = 3156.3333
$ws.Range("K131").Value = 9468.999899999999
$ws.Range("M131").Value = -4428.999899999999

$ws.Range("H132").Value = 874.86365
$ws.Range("I132").Value = 726.0476
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 2178.1428
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = 351.8571999999999
$ws.Range("N132").Value = -17060

$ws.Range("H135").Value = 597.5
$ws.Range("I135").Value = 626.82355
$ws.Range("K135").Value = 5641.41195
$ws.Range("M135").Value = -3106.41195

$ws.Range("H137").Value = 2420.2805
$ws.Range("I137").Value = 1323.5
$ws.Range("K137").Value = 3970.5
$ws.Range("M137").Value = -1420.5

$ws.Range("H138").Value = 2969.4512
$ws.Range("I138").Value = 1670.091
$ws.Range("J138").Value = 3445.8833
$ws.Range("K138").Value = 5010.272999999999
$ws.Range("L138").Value = 10337.6499
$ws.Range("M138").Value = 129.7270000000008
$ws.Range("N138").Value = -20617.6499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4001.6616
$ws.Range("I32").Value = 3324.5
$ws.Range("K32").Value = 3324.5
$ws.Range("M32").Value = -3037.5

$ws.Range("H43").Value = 60170.5
$ws.Range("J43").Value = 63446.668
$ws.Range("L43").Value = 63446.668
$ws.Range("N43").Value = -64072.668

$ws.Range("H74").Value = 16670078
$ws.Range("I74").Value = 27780306
$ws.Range("K74").Value = 27780306
$ws.Range("M74").Value = -27779432

$ws.Range("H77").Value = 16670078
$ws.Range("I77").Value = 27780306
$ws.Range("K77").Value = 138901530
$ws.Range("M77").Value = -138897162

$ws.Range("H97").Value = 2176.6667
$ws.Range("I97").Value = 2176.6667
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 2176.6667
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1680.6667
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3826.6667
$ws.Range("I20").Value = 3471.5789
$ws.Range("J20").Value = 4440
$ws.Range("K20").Value = 3471.5789
$ws.Range("L20").Value = 4440
$ws.Range("M20").Value = -3224.5789
$ws.Range("N20").Value = -4934

$ws.Range("H99").Value = 1462.3334
$ws.Range("I99").Value = 1210
$ws.Range("K99").Value = 1210
$ws.Range("M99").Value = 288

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2467
$ws.Range("I94").Value = 1508.3334
$ws.Range("J94").Value = 3545.5
$ws.Range("K94").Value = 1508.3334
$ws.Range("L94").Value = 3545.5
$ws.Range("M94").Value = -1057.3334
$ws.Range("N94").Value = -4447.5

$ws.Range("H132").Value = 7184.2856
$ws.Range("I132").Value = 5749.5
$ws.Range("J132").Value = 7758.2
$ws.Range("K132").Value = 17248.5
$ws.Range("L132").Value = 23274.6
$ws.Range("M132").Value = -14718.5
$ws.Range("N132").Value = -28334.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 200340.12
$ws.Range("I2").Value = 1001
$ws.Range("J2").Value = 228817.14
$ws.Range("K2").Value = 6006
$ws.Range("L2").Value = 1372902.84
$ws.Range("M2").Value = -5893
$ws.Range("N2").Value = -1373128.84

$ws.Range("H25").Value = 178.16667
$ws.Range("J25").Value = 203.16667
$ws.Range("L25").Value = 609.50001
$ws.Range("N25").Value = -947.50001

$ws.Range("H30").Value = 178.16667
$ws.Range("J30").Value = 203.16667
$ws.Range("L30").Value = 609.50001
$ws.Range("N30").Value = -813.50001

$ws.Range("H131").Value = 7607425
$ws.Range("J131").Value = 5118620
$ws.Range("L131").Value = 15355860
$ws.Range("N131").Value = -15365940

$ws.Range("H136").Value = 2336.5557
$ws.Range("I136").Value = 1575.5714
$ws.Range("K136").Value = 4726.7142
$ws.Range("M136").Value = 373.2857999999997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 34943.332
$ws.Range("J93").Value = 38301
$ws.Range("L93").Value = 38301
$ws.Range("N93").Value = -42045

$ws.Range("H122").Value = 6570.0264
$ws.Range("I122").Value = 5326
$ws.Range("K122").Value = 15978
$ws.Range("M122").Value = -13528

$ws.Range("H132").Value = 3476.0908
$ws.Range("I132").Value = 3023.2727
$ws.Range("K132").Value = 9069.8181
$ws.Range("M132").Value = -6539.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6616.92
$ws.Range("I7").Value = 4771.45
$ws.Range("K7").Value = 4771.45
$ws.Range("M7").Value = -4659.45

$ws.Range("H22").Value = 2773.5
$ws.Range("I22").Value = 1640.65
$ws.Range("K22").Value = 1640.65
$ws.Range("M22").Value = -1345.65

$ws.Range("H27").Value = 2773.5
$ws.Range("I27").Value = 1640.65
$ws.Range("K27").Value = 1640.65
$ws.Range("M27").Value = -1533.65

$ws.Range("H40").Value = 12537.066
$ws.Range("I40").Value = 14721.857
$ws.Range("J40").Value = 10625.375
$ws.Range("K40").Value = 14721.857
$ws.Range("L40").Value = 10625.375
$ws.Range("M40").Value = -14585.857
$ws.Range("N40").Value = -10897.375

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H126").Value = 6616.92
$ws.Range("I126").Value = 4771.45
$ws.Range("K126").Value = 14314.35
$ws.Range("M126").Value = -11844.35

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 11554.333
$ws.Range("J14").Value = 8108.6665
$ws.Range("L14").Value = 8108.6665
$ws.Range("N14").Value = -8444.666499999999

$ws.Range("H56").Value = 48000
$ws.Range("J56").Value = 48000
$ws.Range("L56").Value = 48000
$ws.Range("N56").Value = -49428

$ws.Range("H70").Value = 46403.332
$ws.Range("I70").Value = 50000
$ws.Range("J70").Value = 44605
$ws.Range("K70").Value = 50000
$ws.Range("L70").Value = 44605
$ws.Range("M70").Value = -49685
$ws.Range("N70").Value = -45235

$ws.Range("H73").Value = 46403.332
$ws.Range("I73").Value = 50000
$ws.Range("J73").Value = 44605
$ws.Range("K73").Value = 50000
$ws.Range("L73").Value = 44605
$ws.Range("M73").Value = -48908
$ws.Range("N73").Value = -46789

$ws.Range("H113").Value = 1370.7858
$ws.Range("I113").Value = 998.5
$ws.Range("K113").Value = 2995.5
$ws.Range("M113").Value = -825.5

$ws.Range("H122").Value = 1888.2273
$ws.Range("I122").Value = 1229.8055
$ws.Range("K122").Value = 3689.4165
$ws.Range("M122").Value = -1239.4165

$ws.Range("H126").Value = 1405.2106
$ws.Range("I126").Value = 992
$ws.Range("J126").Value = 2113.5715
$ws.Range("K126").Value = 2976
$ws.Range("L126").Value = 6340.7145
$ws.Range("M126").Value = -506
$ws.Range("N126").Value = -11280.7145

$ws.Range("H132").Value = 10755
$ws.Range("I132").Value = 2505
$ws.Range("J132").Value = 19005
$ws.Range("K132").Value = 7515
$ws.Range("L132").Value = 57015
$ws.Range("M132").Value = -4985
$ws.Range("N132").Value = -62075

